# Generate Report for handback
# For each locale sheet (zh-cn, de-de):
#   - Status (col B) for the two handed-off files flips from
#     "Not yet handed off" to "Handed back"
#   - New "Latest Target File" (col E) / "Latest Handback File" (col F)
#     hyperlinks are added, mirroring the existing "Source File Name" (A)
#     and "Latest Handoff File" (C) links
#   - "Latest Handback DateTime" (col G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

$localeSheets = @(
    @{
        Name = "zh-cn"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8a0ecc7f3dbf5cfbbe69735d66ae96e36cf15dc0/e2e"
        XlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/983dcc705b49df25ecab3efa059f3cd2845e9bc8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
        XlfSuffix = "63d995da3bcf9121499df767357e9ac1b452aa74.zh-cn.xlf"
        XlfSuffix2 = "b40c99edc99d64de8039f4ab4e0228d53227b034.zh-cn.xlf"
        HandbackTime = "2016-01-07 14:15:31"
    },
    @{
        Name = "de-de"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8a0ecc7f3dbf5cfbbe69735d66ae96e36cf15dc0/e2e"
        XlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1ccbf95dd00359c12b7e7ffe8316886aa0c2be7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"
        XlfSuffix = "63d995da3bcf9121499df767357e9ac1b452aa74.de-de.xlf"
        XlfSuffix2 = "b40c99edc99d64de8039f4ab4e0228d53227b034.de-de.xlf"
        HandbackTime = "2016-01-07 14:15:52"
    }
)

$file1 = "1dae8f16-580b-41df-86c4-dea2b1bc2c57"
$file2 = "a6da3a44-0a20-42f7-8585-da3039dcbb2e"

foreach ($loc in $localeSheets) {
    $ws = $wb.Worksheets.Item($loc.Name)

    $md1 = "$file1.md"
    $xlf1 = "$file1.$($loc.XlfSuffix)"
    $md2 = "$file2.md"
    $xlf2 = "$file2.$($loc.XlfSuffix2)"

    # --- Row 2 (1dae8f16...) ---
    $ws.Range("B2").Value = "Handed back"

    $ws.Hyperlinks.Add($ws.Range("E2"), "$($loc.MdUrl)/$md1", "", "", $md1) | Out-Null
    $ws.Range("E2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F2"), "$($loc.XlfBase)/$xlf1", "", "", $xlf1) | Out-Null
    $ws.Range("F2").Style = "HyperLink"

    $ws.Range("G2").Value = $loc.HandbackTime

    # --- Row 3 (a6da3a44...) ---
    $ws.Range("B3").Value = "Handed back"

    $ws.Hyperlinks.Add($ws.Range("E3"), "$($loc.MdUrl)/$md2", "", "", $md2) | Out-Null
    $ws.Range("E3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F3"), "$($loc.XlfBase)/$xlf2", "", "", $xlf2) | Out-Null
    $ws.Range("F3").Style = "HyperLink"

    $ws.Range("G3").Value = $loc.HandbackTime
}

Write-Host "Handback report generated for zh-cn and de-de sheets."
